# Saldo_guide.xlsx update: refresh reference date (IClientBalance extract
# re-run a day later) and correct a handful of balance values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect the newer extraction run.
$ws.Name = "IClientBalance-20241023-095201-"

# Column G ("Dt. Referencia") moves forward one day (2024-10-22 -> 2024-10-23)
# for every data row (2 through 274).
$ws.Range("G2:G274").Value2 = 45588

# A few balances were corrected in this re-run.
$ws.Range("E51").Value2 = 10081.98
$ws.Range("H51").Value2 = 10081.98

$ws.Range("D104").Value2 = -18689.69
$ws.Range("H104").Value2 = 8663.36

$ws.Range("E109").Value2 = 31.98
$ws.Range("H109").Value2 = 31.98
